# Re-apply the "Merge branch 'master' ..." changes that a prior commit had
# reverted: move several task-description cells from columns B/C out to
# columns E/F (freeing up B/C for other content), bump a few row heights to
# match the re-wrapped text, extend row spans to 1:6, and move the active
# selection to B21.

function Move-CellContent($Worksheet, $FromRow, $FromCol, $ToRow, $ToCol) {
    $src = $Worksheet.Cells.Item($FromRow, $FromCol)
    $dst = $Worksheet.Cells.Item($ToRow, $ToCol)
    $src.Copy($dst)
    $src.Clear()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: Passenger -> "Enter Passenger details..." moves C -> F
Move-CellContent $ws 17 3 17 6

# Row 18: Payment (User) -> "Create a payment page..." moves B -> F
Move-CellContent $ws 18 2 18 6
$ws.Rows.Item(18).RowHeight = 86.4

# Row 22: Forgot Password (User) -> "FORGOT PASSWORD SERVICE..." moves C -> F
Move-CellContent $ws 22 3 22 6
$ws.Rows.Item(22).RowHeight = 43.2

# Row 23: Password Validation -> "Strong Password, Confirm Password..." moves C -> F
Move-CellContent $ws 23 3 23 6

# Row 24: Admin View -> "Show all flights..." moves C -> E
Move-CellContent $ws 24 3 24 5
$ws.Rows.Item(24).RowHeight = 34.2

# Row 25: User Profile -> "Show Booking History..." moves B -> E
Move-CellContent $ws 25 2 25 5
$ws.Rows.Item(25).RowHeight = 35.4

# Update the active selection to B21 (was G14)
$ws.Range("B21").Select()
